# Regenerate orders with updated distance/sizes.
# Performs a global token substitution across the Condition, Filename_Left,
# Filename_Right, Distance and Size columns:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Columns that contain text tokens needing replacement: B, D, E, H, J
$targetCols = @(2, 4, 5, 8, 10)

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()

        if ($val -ne $null -and $val -is [string]) {
            $newVal = $val
            $newVal = $newVal.Replace("D64", "D69")
            $newVal = $newVal.Replace("D51", "D55")
            $newVal = $newVal.Replace("D80", "D86")
            $newVal = $newVal.Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
